$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 9; $row++) {
        $ws.Cells.Item($row, 6).Value = 0
    }
}
